$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New skill rows (10-16), matching column layout:
# A=name B=description C=str_mod D=dex_mod E=dur_mod F=chr_mod G=focus_mod
# H=int_mod I=agi_mod J=base_damage_mod K=base_ac_mod L=base_healing_mod
# M=max_level N=total_kills_needed O=parent_id P=parent_level_needed

# Row 10 - Clerics Prayer
$ws.Range("A10").Value = "Clerics Prayer"
$ws.Range("B10").Value = "Pray to your god child, pray for the light to come and heal your wounds. While training this skill, over time you will gain 2% towards your CHR for a total of 10% You will also gain 3% towards healing for an extra 15%"
$ws.Range("F10").Value = 0.02
$ws.Range("L10").Value = 0.03
$ws.Range("M10").Value = 5
$ws.Range("N10").Value = 175

# Row 11 - Clerics of War
$ws.Range("A11").Value = "Clerics of War"
$ws.Range("B11").Value = "Death comes to those who are not prepared. Raise your ac and healing by 5% per level for a total of 15% at max level"
$ws.Range("K11").Value = 0.05
$ws.Range("L11").Value = 0.05
$ws.Range("M11").Value = 3
$ws.Range("N11").Value = 300
$ws.Range("O11").Value = "Clerics Prayer"
$ws.Range("P11").Value = 3

# Row 12 - Durable Priest
$ws.Range("A12").Value = "Durable Priest"
$ws.Range("B12").Value = "Stand fast good sir! Over time raise your defense and durability by 4% for a total of 20% at max level."
$ws.Range("E12").Value = 0.04
$ws.Range("K12").Value = 0.04
$ws.Range("M12").Value = 5
$ws.Range("N12").Value = 600
$ws.Range("O12").Value = "Clerics Prayer"
$ws.Range("P12").Value = 3

# Row 13 - Clerics Wrath
$ws.Range("A13").Value = "Clerics Wrath"
$ws.Range("B13").Value = "Lash out child. Lash out with the words of your god. Raise your focus and damage over time by 5% for a total of 25% at max level."
$ws.Range("G13").Value = 0.05
$ws.Range("J13").Value = 0.05
$ws.Range("M13").Value = 5
$ws.Range("N13").Value = 1000
$ws.Range("O13").Value = "Clerics Prayer"
$ws.Range("P13").Value = 4

# Row 14 - Prophets Grace
$ws.Range("A14").Value = "Prophets Grace"
$ws.Range("B14").Value = "Stand in the field of battle and protect those around you and heal your wounds. Over time you will gain an additional 40% towards your Attack, 50% towards your Armour Class and 75% towards your healing."
$ws.Range("J14").Value = 0.08
$ws.Range("K14").Value = 0.1
$ws.Range("L14").Value = 0.15
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 1000
$ws.Range("O14").Value = "Clerics of War"
$ws.Range("P14").Value = 2

# Row 15 - Prophets Rage
$ws.Range("A15").Value = "Prophets Rage"
$ws.Range("B15").Value = "Rage at the enemy but in a godly way. Raise your Attack and CHR by 40% and 100% at max level."
$ws.Range("F15").Value = 0.25
$ws.Range("J15").Value = 0.1
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = 1200
$ws.Range("O15").Value = "Clerics Wrath"
$ws.Range("P15").Value = 3

# Row 16 - Godly Cosmic Awakening
$ws.Range("A16").Value = "Godly Cosmic Awakening"
$ws.Range("B16").Value = "Awaken the power of the one true god according to The Churches doctrine."
$ws.Range("F16").Value = 0.25
$ws.Range("G16").Value = 0.25
$ws.Range("J16").Value = 0.25
$ws.Range("K16").Value = 0.25
$ws.Range("L16").Value = 0.25
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 1500
$ws.Range("O16").Value = "Prophets Rage"
$ws.Range("P16").Value = 2

# Columns A and B grow to fit the new, longer skill names/descriptions
# (matches Excel's pixel-grid column-width quantization for the target widths)
$ws.Range("A1").ColumnWidth = 26.16
$ws.Range("B1").ColumnWidth = 256.16

